# Add "NA" under duplicate_image_filename (column E) for the practice
# and main trial rows (rows 2-21) of the stimuli sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E21").Value = "NA"
